# Atualização Visual e Funcional do programa
# Atualiza o status exibido e a forma de identificação do disco do computador.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Status: "Ok" -> "Ativo"
$ws.Range("D2").Value = "Ativo"

# HD: "HDD - 931,5 GB" -> "SSD 931,5 GB"
$ws.Range("K2").Value = "SSD 931,5 GB"
